# Update Name of Algo - apply corrected imputed values to result_data_RandomForest sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = -6.73279999999999
$ws.Range("D3").Value  = -7.585299999999992
$ws.Range("D5").Value  = -8.080499999999999
$ws.Range("C9").Value  = -11.89620000000001
$ws.Range("D11").Value = -8.368100000000004
$ws.Range("D12").Value = -8.367000000000006
$ws.Range("C13").Value = -12.91409999999999
$ws.Range("C16").Value = -11.6239
$ws.Range("C18").Value = -14.24049999999998
$ws.Range("C20").Value = -13.83369999999998
$ws.Range("D21").Value = -7.6842
